# Weekly data update: a new price record (week of 2021-10-15, date serial 44484)
# is inserted for "Femacal de La Calera" / "Zapallo italiano" / "Primera",
# pushing every subsequent row down by one and extending the sheet from
# A1:R226 to A1:R227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 117 (row 117 and everything below it
# shifts down by one; the former last row, 226, becomes row 227).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A117").Value = 3
$ws.Range("B117").Value = "Femacal de La Calera"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44484
$ws.Range("E117").Value = 5
$ws.Range("F117").Value = 100112032
$ws.Range("G117").Value = "Zapallo italiano"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 190
$ws.Range("K117").Value = 11000
$ws.Range("L117").Value = 12500
$ws.Range("M117").Value = 11763
$ws.Range("N117").Value = "`$/caja 70 unidades"
$ws.Range("O117").Value = "Región de Arica y Parinacota"
$ws.Range("P117").Value = 168
$ws.Range("Q117").Value = 70
$ws.Range("R117").Value = "Hortaliza"
